$d = $word.ActiveDocument
$r = $d.Content
$r.Collapse(0)

$xmlFragment = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:cx="http://schemas.microsoft.com/office/drawing/2014/chartex" xmlns:cx1="http://schemas.microsoft.com/office/drawing/2015/9/8/chartex" xmlns:cx2="http://schemas.microsoft.com/office/drawing/2015/10/21/chartex" xmlns:cx3="http://schemas.microsoft.com/office/drawing/2016/5/9/chartex" xmlns:cx4="http://schemas.microsoft.com/office/drawing/2016/5/10/chartex" xmlns:cx5="http://schemas.microsoft.com/office/drawing/2016/5/11/chartex" xmlns:cx6="http://schemas.microsoft.com/office/drawing/2016/5/12/chartex" xmlns:cx7="http://schemas.microsoft.com/office/drawing/2016/5/13/chartex" xmlns:cx8="http://schemas.microsoft.com/office/drawing/2016/5/14/chartex" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:aink="http://schemas.microsoft.com/office/drawing/2016/ink" xmlns:am3d="http://schemas.microsoft.com/office/drawing/2017/model3d" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:oel="http://schemas.microsoft.com/office/2019/extlst" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:w16cex="http://schemas.microsoft.com/office/word/2018/wordml/cex" xmlns:w16cid="http://schemas.microsoft.com/office/word/2016/wordml/cid" xmlns:w16="http://schemas.microsoft.com/office/word/2018/wordml" xmlns:w16sdtdh="http://schemas.microsoft.com/office/word/2020/wordml/sdtdatahash" xmlns:w16se="http://schemas.microsoft.com/office/word/2015/wordml/symex" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 w16se w16cid w16 w16cex w16sdtdh wp14">
<w:body>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListeParagraf"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="4082"/>
          <w:tab w:val="left" w:pos="7864"/>
        </w:tabs>
        <w:ind w:left="1080"/>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListeParagraf"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="4082"/>
          <w:tab w:val="left" w:pos="7864"/>
        </w:tabs>
        <w:ind w:left="1080"/>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListeParagraf"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="4082"/>
          <w:tab w:val="left" w:pos="7864"/>
        </w:tabs>
        <w:ind w:left="1080"/>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListeParagraf"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="4082"/>
          <w:tab w:val="left" w:pos="7864"/>
        </w:tabs>
        <w:ind w:left="1080"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>Sanal Bilgisayar Kurulumları ve Ayarları</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListeParagraf"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="4082"/>
          <w:tab w:val="left" w:pos="7864"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
        <w:t>Virtual Box indir ve kurulum adımlarına başla</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListeParagraf"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="4082"/>
          <w:tab w:val="left" w:pos="7864"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
        <w:t>Eğer ki kurulum esnasında hata verir ise Visual Studio C++ v_credit dosyası Windows’un kendi sitesinden indilir ve next -&gt; next  diyerek kurulum tamamlanır.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListeParagraf"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="4082"/>
          <w:tab w:val="left" w:pos="7864"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
        <w:t>Ayrıca C++ dosyası kurulduktan sonra da virtual box programı next -&gt; next diyerek kurulumu tamamlanır.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListeParagraf"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="4082"/>
          <w:tab w:val="left" w:pos="7864"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
        <w:t>ISO Dosyası Yükleme : Katılımsız kurulumu atla seçeneği seçilir ve gerekli olan ram miktarı belirlenir. Daha sonra ayrılacak olan disk miktarı seçilir. (genelde varsayılan olarak gelen değer bırakılır) Tamam denilerek işlem tamamlanır.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListeParagraf"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="4082"/>
          <w:tab w:val="left" w:pos="7864"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
        <w:t>Sistemi başlatmadan önce ayarlardan ram miktarları değiştirilebilir.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListeParagraf"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="4082"/>
          <w:tab w:val="left" w:pos="7864"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
        <w:t>Standart Windows kurulum aşamalarına devam edilir. Kurulum esnasında daha yönetilebilir bir sistem için Windows 10 Pro seçilir.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListeParagraf"/>
        <w:numPr>
          <w:ilvl w:val="1"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="4082"/>
          <w:tab w:val="left" w:pos="7864"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
        <w:t>Kurulumdan sonra 2 ayar yapılmalıdır.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListeParagraf"/>
        <w:numPr>
          <w:ilvl w:val="2"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="4082"/>
          <w:tab w:val="left" w:pos="7864"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
        <w:t>Settings -&gt; storage -&gt; .iso tıklanır -&gt; sağ kısımdan remove disk</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListeParagraf"/>
        <w:numPr>
          <w:ilvl w:val="2"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="4082"/>
          <w:tab w:val="left" w:pos="7864"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
        <w:t>Settings -&gt; network -&gt; bridge mode yapılır.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListeParagraf"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="4082"/>
          <w:tab w:val="left" w:pos="7864"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
        <w:t xml:space="preserve">Otomatik ekran genişliği : Açık olan pencerede insert guest additions tıklanır. CD-Rom sürücüsüne eklenmiş olan  uygulama next -&gt; next denilerek </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
        <w:t>kuruluma denilerek kuruluma devam edilir. Daha sonra sistem yeniden başlatılınca view -&gt; guest auto resize display seçeneği seçilir. Daha sonra oluşan .iso dosyası storagedan kaldırılır.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListeParagraf"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="4082"/>
          <w:tab w:val="left" w:pos="7864"/>
        </w:tabs>
        <w:ind w:left="1080"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
        <w:t>İki bilgisayarın haberleşmesi</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListeParagraf"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="4082"/>
          <w:tab w:val="left" w:pos="7864"/>
        </w:tabs>
        <w:ind w:left="1080"/>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
        <w:t>IP atama : sağ alt kısımdan monitör -&gt; ağ ve internet ayarları -&gt; ethernet -&gt; bağdaştırıcı seçenekleri -&gt; çift tıklama -&gt; IPv4 -&gt; çift tıklama -&gt; aşağıdaki IP adresini kullan seçeneğine gerekli IP bilgileri girilir. Tab tuşuna basıldıktan sonra pc tarafından girilen IP bilgisine göre ağ alt maskesi otomatik olarak verilir.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListeParagraf"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="4082"/>
          <w:tab w:val="left" w:pos="7864"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
        <w:t>Otomatik IP atamalarında cihaz IP alamazsa 168.xxx.xxx.xxx şeklinde bir IP numarası görünür ise APIPA’ya düşmüş demektir.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListeParagraf"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="4082"/>
          <w:tab w:val="left" w:pos="7864"/>
        </w:tabs>
        <w:ind w:left="1080"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
        <w:t>Bağlantı Testleri (ping) ve IP öğrenme</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListeParagraf"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="4082"/>
          <w:tab w:val="left" w:pos="7864"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
        <w:t>CMD -&gt; ping 192.168.1.xxx -&gt; enter yapılınca eğer başarılı bir şekilde geri dönüş sağlanıyorsa bağlantı başarılıdır.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListeParagraf"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="4082"/>
          <w:tab w:val="left" w:pos="7864"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
        <w:t>IP Öğrenme : CMD -&gt; ipconfig -&gt; enter tıklanınca IP bilgisini gösterir.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListeParagraf"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="4082"/>
          <w:tab w:val="left" w:pos="7864"/>
        </w:tabs>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS"/>
        </w:rPr>
        <w:t>Ping -t -&gt; normal ping komutu 4 satırlık bir geri dönüş sağlar. Eğer ping -t IP_bilgisi şeklinde bir komut verilirse sonsuza kadar istek atmaya devam eder.</w:t>
       </w:r>
     </w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@

$r.InsertXML($xmlFragment)
Write-Output "Inserted new section content."
